# The presentation's slide master currently uses the "Integral" theme
# (ppt/theme/theme1.xml). This edit re-colors that theme to match the
# stock "Office Theme" color scheme (the swap recorded in the target
# diff), by rewriting each of the twelve theme colors via the
# PowerPoint object model's ThemeColorScheme on the master theme.
#
# MsoThemeColorSchemeIndex order used by ThemeColorScheme.Item():
#   1 = dk1 (Dark 1)        7  = accent3
#   2 = lt1 (Light 1)       8  = accent4
#   3 = dk2 (Dark 2)        9  = accent5
#   4 = lt2 (Light 2)       10 = accent6
#   5 = accent1             11 = hlink
#   6 = accent2              12 = folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
